# Add the Arabic language row to the registration-center-type master data sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row 4 values ---
$ws.Range("A4").Value = "ara"
$ws.Range("B4").Value = "REG"
$ws.Range("C4").Value = "عادي"
$ws.Range("D4").Value = "مركز التسجيل العادي"

# Keep "TRUE" as literal text (shared string), matching E2/E3, not a boolean.
$e4 = $ws.Range("E4")
$e4.Formula = '="TRU"&"E"'
$e4.Copy()
$e4.PasteSpecial(-4163)
$excel.CutCopyMode = $false

# --- Row height for the new row ---
$ws.Rows.Item(4).RowHeight = 16.4

# --- Column widths for C and D (new, wider columns to fit the Arabic text) ---
$ws.Columns.Item(3).ColumnWidth = 21.1
$ws.Columns.Item(4).ColumnWidth = 33.1

# --- Cell formatting for C4/D4: left aligned, wrap text ---
$c4 = $ws.Range("C4")
$c4.HorizontalAlignment = -4131
$c4.WrapText = $true

$c4.Copy()
$ws.Range("D4").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Update selection to match the newly-entered cells ---
$ws.Range("C4:D4").Select() | Out-Null
